$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape "文本框 3" (id=4): speech bubble "kyon:「What should I do ..."
# Split the leading speaker name into its own run: "kyon" -> "Kyon"
# ---------------------------------------------------------------------------
$sh4 = $s.Shapes.Item(4)
$tr4 = $sh4.TextFrame.TextRange
$tr4.Characters(1, 4).Text = "Kyon"
# this shape auto-fits its box to the text; re-editing the run recomputes
# the fitted height via font metrics, so pin it back to the stored value.
$sh4.Height = 21.7

# ---------------------------------------------------------------------------
# Shape "文本框 5" (id=6): speech bubble "haruhi:「You're right, ..."
# Split the leading speaker name into its own run and give it its own
# trailing colon run as well: "haruhi:" -> "Haruhi" + ":"
# ---------------------------------------------------------------------------
$sh5 = $s.Shapes.Item(5)
$tr5 = $sh5.TextFrame.TextRange
$tr5.Characters(7, 1).Text = ":"
$tr5.Characters(1, 6).Text = "Haruhi"

# ---------------------------------------------------------------------------
# Shape "文本框 6" (id=7): speech bubble "kyon:「If we form a club, ..."
# "kyon:" -> "Kyon" + ":"
# ---------------------------------------------------------------------------
$sh6 = $s.Shapes.Item(6)
$tr6 = $sh6.TextFrame.TextRange
$tr6.Characters(5, 1).Text = ":"
$tr6.Characters(1, 4).Text = "Kyon"

# ---------------------------------------------------------------------------
# Shape "文本框 8" (id=9): speech bubble with the club-naming monologue.
# - reposition the box
# - "...haruhi:" -> "...Haruhi" + ":" (and drop one of the leading spaces)
# - "Haruhi Suzumiya" -> "Haruhi" + " " + "Suzumiya" (own runs)
# ---------------------------------------------------------------------------
$sh9 = $s.Shapes.Item(7)
$tr9 = $sh9.TextFrame.TextRange

# split "Haruhi" / " " / "Suzumiya" in the second paragraph first (highest
# offsets), then the name near the start, then drop a leading space -- all
# done high-offset-first so earlier offsets stay valid while we work.
$tr9.Characters(155, 8).Text = "Suzumiya"
$tr9.Characters(148, 6).Text = "Haruhi"
$tr9.Characters(50, 1).Text = ":"
$tr9.Characters(44, 6).Text = "Haruhi"
$tr9.Characters(43, 1).Text = ""

# Move the shape from (-67945, 3188970) EMU to (82884, 3174047) EMU.
# Shape.Left/Top are single-precision (points); pick values whose f32
# round-trip lands exactly on the target EMU.
$sh9.Left = 6.526299312598425
$sh9.Top = 249.92496492992126
